# Harp expander BOM update: add clock generator parts (IC8, IC9, IC10, C35, C36, R54, R55),
# move R53 from the "1k" group to the "10R" group, and remove the obsolete
# "CLCK Jumper" (JP1/JP2) BOM line since the clock jack is now used directly
# as a clock generator input/output (no jumper needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Qty / part-list updates on the existing rows (before the row delete,
#        while row numbers still match the original layout) -----------------

# Row 3: 100nF / C0402 -> qty 24 -> 26, add C35, C36
$ws.Range("A3").Value = 26
$ws.Range("D3").Value = "C1, C3, C5, C7, C9, C12, C13, C14, C15, C16, C17, C18, C19, C20, C21, C22, C25, C26, C27, C28, C29, C30, C31, C32, C35, C36"

# Row 4: 10R / R0402 -> qty 1 -> 2, add R53 (moved from the "1k" group)
$ws.Range("A4").Value = 2
$ws.Range("D4").Value = "R2, R53"

# Row 5: 10k / R0402 -> qty 16 -> 18, add R54, R55
$ws.Range("A5").Value = 18
$ws.Range("D5").Value = "R3, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19, R38, R43, R46, R47, R54, R55"

# Row 7: 1k / R0402 -> qty 13 -> 12, drop R53 (moved to the "10R" group)
$ws.Range("A7").Value = 12
$ws.Range("D7").Value = "R20, R21, R22, R23, R24, R25, R26, R27, R28, R29, R30, R31"

# Row 25: SN74LVC1G125DBVT -> qty 2 -> 4, IC7/IC14 -> IC7/IC8/IC9/IC10
$ws.Range("A25").Value = 4
$ws.Range("D25").Value = "IC7, IC8, IC9, IC10"

# --- 2. Remove the obsolete "CLCK Jumper" BOM row (row 18) ------------------
# Deleting the row shifts every following row up by one automatically.
$ws.Rows.Item(18).Delete()

# --- 3. Fix up the summary block at the bottom (row numbers already
#        shifted up by one after the delete above) ---------------------------

# Clear the now-orphaned helper formula that used to live in G42 (shifted to G41)
$ws.Range("G41").ClearContents()

# Number of unique parts / SMD parts / TH parts counts
$ws.Range("F40").Value = 27
$ws.Range("F41").Value = 122
$ws.Range("F42").Value = 13

# --- 4. Sheet view tidy-up (matches the author's re-saved view state) -------
$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 1
